# Updates the "cryptos" price/volume table with refreshed figures, and
# swaps the Cosmos/Monero rows (25 <-> 26) to match the new ranking order.
# Price-looking values are written via NumberFormat="@" + Style="Normal"
# so they land as plain text (matching the source data, which stores
# prices as text because "." is used both as thousands & decimal sep.)
# without leaving a residual cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = '30.589.97'
$ws.Cells.Item(2, 5).Formula = '  -0.59%  '
$ws.Cells.Item(3, 4).Formula = '1.875.08'
$ws.Cells.Item(3, 5).Formula = '  -0.94%  '
$ws.Cells.Item(4, 5).Formula = '  -0.07%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Formula = '247.89'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Formula = '  +1.04%  '
$ws.Cells.Item(6, 5).Formula = '  -0.02%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Formula = '0.4753'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Formula = '  -0.69%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Formula = '0.2907'
$c.Style = 'Normal'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Formula = '0.06486'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Formula = '  -1.23%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Formula = '21.97'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Formula = '  +2.79%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Formula = '0.07750'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Formula = '  -0.40%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Formula = '0.7379'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Formula = '  -1.17%  '
$ws.Cells.Item(13, 4).Formula = '1.875.31'
$ws.Cells.Item(13, 5).Formula = '  -0.95%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Formula = '96.04'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Formula = '  -1.04%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Formula = '5.183'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Formula = '  -0.07%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Formula = '274.15'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Formula = '  -2.41%  '
$ws.Cells.Item(17, 4).Formula = '30.631.75'
$ws.Cells.Item(17, 5).Formula = '  -0.46%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Formula = '13.23'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Formula = '  -2.08%  '
$ws.Cells.Item(19, 5).Formula = '  -0.03%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Formula = '0.000007492'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Formula = '  -1.67%  '
$ws.Cells.Item(21, 4).Formula = '2.121.32'
$ws.Cells.Item(21, 5).Formula = '  -1.50%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Formula = '0.9993'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Formula = '  -0.22%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Formula = '5.221'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Formula = '  -1.57%  '
$ws.Cells.Item(24, 5).Formula = '  -1.13%  '
$ws.Cells.Item(25, 2).Formula = 'Cosmos'
$ws.Cells.Item(25, 3).Formula = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Formula = '9.199'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Formula = '  -1.71%  '
$ws.Cells.Item(26, 2).Formula = 'Monero'
$ws.Cells.Item(26, 3).Formula = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Formula = '165.08'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Formula = '  -0.68%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Formula = '18.79'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Formula = '  -1.87%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Formula = '1.907'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Formula = '  -3.46%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Formula = '0.09892'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Formula = '  -1.14%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Formula = '1.345'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Formula = '  -1.97%  '
$ws.Cells.Item(31, 5).Formula = '  -0.70%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Formula = '4.259'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Formula = '  -2.67%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Formula = '4.095'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Formula = '  -0.97%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Formula = '0.04781'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Formula = '  -0.28%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Formula = '1.120'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Formula = '  -1.04%  '
$ws.Cells.Item(36, 5).Formula = '  -1.55%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Formula = '2.719'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Formula = '  +0.02%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Formula = '0.01852'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Formula = '  -1.47%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Formula = '2.762'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Formula = '  -0.30%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Formula = '6.265'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Formula = '  -2.62%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Formula = '73.42'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Formula = '  +3.96%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Formula = '1.981'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Formula = '  +2.53%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Formula = '0.4181'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Formula = '  -0.98%  '
$ws.Cells.Item(44, 5).Formula = '  -0.03%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Formula = '0.8356'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Formula = '  -1.69%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Formula = '101.60'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Formula = '  -1.05%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Formula = '9.367'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Formula = '  -0.82%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Formula = '35.41'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Formula = '  +0.16%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Formula = '6.971'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Formula = '  -2.93%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Formula = '918.61'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Formula = '  -2.58%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Formula = '0.05670'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Formula = '  +0.85%  '
